# Apply the "addind preparer to sheet" change:
#  - set E2:E19 (purpose column) to "fullRNASEQ" instead of "S.GISH"
#  - update the sheet's selection to D20:F24 (active cell D20)
#  - enable iterative calculation with a max change (iterateDelta) of 1E-4

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update purpose column values for the data rows
$ws.Range("E2:E19").Value = "fullRNASEQ"

# Move the selection as recorded in the saved workbook
$ws.Range("D20:F24").Select()

# Turn on iterative calculation with the new max-change delta
$excel.Iteration = $true
$excel.MaxChange = 0.0001
